$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 5 new "batch 32" device rows (157-161), mirroring the existing pattern
# of Finger Print Scanner / IRIS Scanner / Web Camera / Document Scanner /
# Printer entries used for every previous batch in this sheet.
# ---------------------------------------------------------------------------

$names   = @("Finger Print Scanner 32", "IRIS Scanner 32", "Web Camera 32", "Document Scanner 32", "Printer 32")
$macs    = @("80-75-40-E8-CA-24", "0E-1A-14-4A-6D-3A", "65-13-7F-0F-F7-53", "73-C4-DE-8E-C9-8D", "EC-74-AB-E0-0F-38")
$serials = @("BS563Q2230824", "BS563Q2230825", "BS563Q2230826", "BS563Q2230827", "BS563Q2230828")
$dspecs  = @(165, 327, 736, 801, 920)

$firstRow = 157
$lastDataRow = 156

# Write column-by-column (not row-by-row) so new shared-string entries land
# in the same order as the reference workbook: all names, then all MAC
# addresses, then all serial numbers.
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = 3000176 + $i
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($firstRow + $i, 2).Value = $names[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($firstRow + $i, 3).Value = $macs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($firstRow + $i, 4).Value = $serials[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($firstRow + $i, 6).Value = $dspecs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($lastDataRow, 7).Value()
    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($lastDataRow, 9).Value()
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($lastDataRow, 10).Value()
}

# ---------------------------------------------------------------------------
# Update the saved view state: scrolled/selected further down than before.
# ---------------------------------------------------------------------------
$ws.Range($ws.Columns.Item(11), $ws.Columns.Item(16384)).Select() | Out-Null
